$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Disease Ontology source_version (cell E3): v2024-02-28 -> v2024-03-28
$ws.Range("E3").Value = "v2024-03-28"

# Move/restore the cursor selection to E3 (matches the saved cursor position)
$ws.Range("E3").Select()
